# Corporate Customer excel file update:
# Add new header columns (G1:L1) for manager/locker details, refresh the
# dimension/selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells appended after the existing F1 (OPEARTING.NAME:1) header.
$ws.Range("G1").Value = "OLD.LOCKER.NO"
$ws.Range("H1").Value = "MAN.NAME:1"
$ws.Range("I1").Value = "MAN.ADDRESS:1:1"
$ws.Range("J1").Value = "MAN.CONTACT.NO:1"
$ws.Range("K1").Value = "MAN.ID.TYPE:1"
$ws.Range("L1").Value = "MAN.ID.NO:1"

# Matches the selection recorded in the saved workbook (active cell moved to F12).
$ws.Range("F12").Select()
